$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Beggs & Graddy (2009) argue that demand for art changes," becomes
#    "... demand for art changes over time,"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "demand for art changes,", $true, $false, $false, $false, $false,
    $true, 1, $false, "demand for art changes over time,", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "... underlying hedonic quality is constant across sales, allowing for
#    anchoring to be isolated." becomes "... quality remains constant, thus
#    allowing anchoring to be identified."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "quality is constant across sales, allowing for anchoring to be isolated",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "quality remains constant, thus allowing anchoring to be identified", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "... given we control for those quality changes." becomes
#    "... given we control for those differences."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "given we control for those quality changes.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "given we control for those differences.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) New sentence opener inserted before "The mechanism by which ..." and
#    "is a black box" gains "still".
#    "The mechanism by which past quantities impact future ones is a black
#    box," becomes "Even after controlling for such factors, the mechanism
#    by which past quantities impact future ones is still a black box,"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The mechanism by which past quantities impact future ones is a black box,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Even after controlling for such factors, the mechanism by which past quantities impact future ones is still a black box,",
    2) | Out-Null

# The "_GoBack" bookmark (Word's "last edit location" marker) now belongs
# right after the newly-typed "Even after", rather than its old spot near
# "Below". Move it there.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$t = $d.Content.Text
$marker = "Even after"
$idx = $t.IndexOf($marker)
$posAfter = $idx + $marker.Length
$r = $d.Range($posAfter, $posAfter)
$r.Bookmarks.Add("_GoBack") | Out-Null
